# Regenerate save_data: column G ("K") is recalculated from the new
# per-game strikeout source (replacing the old "Strike#" pitch-count
# derived values) and written back as literal values (s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K value (only rows whose K value actually changed
# are listed; rows 41 and 42 keep their original K of 0).
$kValues = [ordered]@{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 2
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 2
    18 = 3
    19 = 2
    20 = 1
    21 = 1
    22 = 0
    23 = 3
    24 = 1
    25 = 1
    26 = 1
    27 = 0
    28 = 1
    29 = 1
    30 = 0
    31 = 2
    32 = 2
    33 = 1
    34 = 1
    35 = 2
    36 = 1
    37 = 1
    38 = 1
    39 = 0
    40 = 1
    43 = 2
    44 = 2
    45 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
